# Update Sheets via scheduled runner: refresh market-price derived values
$wb = $excel.ActiveWorkbook

# Sheet ALC, row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 39838.86
$ws.Range("J17").Value = 40541.34
$ws.Range("L17").Value = 121624.02
$ws.Range("N17").Value = -121960.02

# Sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3208.327
$ws.Range("I137").Value = 3262.923
$ws.Range("J137").Value = 3044.5386
$ws.Range("K137").Value = 9788.769
$ws.Range("L137").Value = 9133.6158
$ws.Range("M137").Value = -7238.769
$ws.Range("N137").Value = -14233.6158

# Sheet ARM, row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7260.2524
$ws.Range("I32").Value = 5991.978
$ws.Range("J32").Value = 21686.875
$ws.Range("K32").Value = 5991.978
$ws.Range("L32").Value = 21686.875
$ws.Range("M32").Value = -5704.978
$ws.Range("N32").Value = -22260.875

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1983.2609
$ws.Range("I61").Value = 1539.2903
$ws.Range("K61").Value = 1539.2903
$ws.Range("M61").Value = -1327.2903

# Sheet ARM, row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1219.2609
$ws.Range("I74").Value = 853.1053000000001
$ws.Range("J74").Value = 2958.5
$ws.Range("K74").Value = 853.1053000000001
$ws.Range("L74").Value = 2958.5
$ws.Range("M74").Value = 20.89469999999994
$ws.Range("N74").Value = -4706.5

# Sheet ARM, row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1219.2609
$ws.Range("I77").Value = 853.1053000000001
$ws.Range("J77").Value = 2958.5
$ws.Range("K77").Value = 4265.5265
$ws.Range("L77").Value = 14792.5
$ws.Range("M77").Value = 102.4735000000001
$ws.Range("N77").Value = -23528.5

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2694.2273
$ws.Range("I132").Value = 2063.8064
$ws.Range("J132").Value = 4197.5386
$ws.Range("K132").Value = 6191.4192
$ws.Range("L132").Value = 12592.6158
$ws.Range("M132").Value = -3661.4192
$ws.Range("N132").Value = -17652.6158

# Sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1983.2609
$ws.Range("I136").Value = 1539.2903
$ws.Range("K136").Value = 4617.8709
$ws.Range("M136").Value = -2067.8709

# Sheet BSM, row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 665.93335
$ws.Range("I94").Value = 578.9
$ws.Range("J94").Value = 840
$ws.Range("K94").Value = 578.9
$ws.Range("L94").Value = 840
$ws.Range("M94").Value = -127.9
$ws.Range("N94").Value = -1742

# Sheet BSM, row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2996.7292
$ws.Range("I134").Value = 3158.5757
$ws.Range("J134").Value = 2640.6667
$ws.Range("K134").Value = 9475.7271
$ws.Range("L134").Value = 7922.000100000001
$ws.Range("M134").Value = -6940.7271
$ws.Range("N134").Value = -12992.0001

# Sheet CRP, row 19
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 138.57143
$ws.Range("I19").Value = 145
$ws.Range("J19").Value = 100
$ws.Range("K19").Value = 145
$ws.Range("L19").Value = 100
$ws.Range("M19").Value = 25
$ws.Range("N19").Value = -440

# Sheet CRP, row 24
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H24").Value = 138.57143
$ws.Range("I24").Value = 145
$ws.Range("J24").Value = 100
$ws.Range("K24").Value = 145
$ws.Range("L24").Value = 100
$ws.Range("M24").Value = 25
$ws.Range("N24").Value = -440

# Sheet CRP, row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2667.8965
$ws.Range("I132").Value = 2797.9285
$ws.Range("J132").Value = 2546.5334
$ws.Range("K132").Value = 8393.7855
$ws.Range("L132").Value = 7639.600199999999
$ws.Range("M132").Value = -5863.7855
$ws.Range("N132").Value = -12699.6002

# Sheet CRP, row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 17245814
$ws.Range("I134").Value = 23814414
$ws.Range("J134").Value = 3239.25
$ws.Range("K134").Value = 71443242
$ws.Range("L134").Value = 9717.75
$ws.Range("M134").Value = -71440707
$ws.Range("N134").Value = -14787.75

# Sheet CUL, row 86
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 1980
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 1980
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 5940
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -8312

# Sheet CUL, row 89
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 1980
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 1980
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 17820
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -29676

# Sheet CUL, row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1272.6666
$ws.Range("I131").Value = 3115
$ws.Range("J131").Value = 1052.6865
$ws.Range("K131").Value = 9345
$ws.Range("L131").Value = 3158.0595
$ws.Range("M131").Value = -4305
$ws.Range("N131").Value = -13238.0595

# Sheet GSM, row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3022.6545
$ws.Range("I132").Value = 2771.8057
$ws.Range("J132").Value = 3497.9473
$ws.Range("K132").Value = 8315.417099999999
$ws.Range("L132").Value = 10493.8419
$ws.Range("M132").Value = -5785.417099999999
$ws.Range("N132").Value = -15553.8419

# Sheet LTW, row 96
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H96").Value = 30000
$ws.Range("J96").Value = 30000
$ws.Range("L96").Value = 30000
$ws.Range("N96").Value = -35492

# Sheet LTW, row 97
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H97").Value = 27333
$ws.Range("J97").Value = 27333
$ws.Range("L97").Value = 27333
$ws.Range("N97").Value = -29315

# Sheet LTW, row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3959.034
$ws.Range("I132").Value = 1571.7916
$ws.Range("J132").Value = 5596
$ws.Range("K132").Value = 4715.3748
$ws.Range("L132").Value = 16788
$ws.Range("M132").Value = -2185.3748
$ws.Range("N132").Value = -21848

# Sheet LTW, row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2515.2327
$ws.Range("I136").Value = 2474.3103
$ws.Range("J136").Value = 2600
$ws.Range("K136").Value = 7422.9309
$ws.Range("L136").Value = 7800
$ws.Range("M136").Value = -4872.9309
$ws.Range("N136").Value = -12900

# Sheet LTW, row 140
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 64666.668
$ws.Range("J140").Value = 64666.668
$ws.Range("L140").Value = 64666.668
$ws.Range("N140").Value = -75026.66800000001

# Sheet WVR, row 97
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 19333.334
$ws.Range("I97").Value = 23000
$ws.Range("J97").Value = 18600
$ws.Range("K97").Value = 23000
$ws.Range("L97").Value = 18600
$ws.Range("M97").Value = -22009
$ws.Range("N97").Value = -20582

# Sheet WVR, row 99
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H99").Value = 27600
$ws.Range("I99").Value = 18000
$ws.Range("K99").Value = 18000
$ws.Range("M99").Value = -15005

# Sheet WVR, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 21587.357
$ws.Range("I132").Value = 3557.6667
$ws.Range("J132").Value = 54040.8
$ws.Range("K132").Value = 10673.0001
$ws.Range("L132").Value = 162122.4
$ws.Range("M132").Value = -8143.000100000001
$ws.Range("N132").Value = -167182.4

# Sheet WVR, row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3717.9443
$ws.Range("I136").Value = 3411.4546
$ws.Range("J136").Value = 4199.5713
$ws.Range("K136").Value = 10234.3638
$ws.Range("L136").Value = 12598.7139
$ws.Range("M136").Value = -7684.363799999999
$ws.Range("N136").Value = -17698.7139
